$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "301.40"),
    @("E2", "-0.73%"),
    @("E3", "-1.88%"),
    @("D4", "5.111"),
    @("E4", "-2.50%"),
    @("D5", "0.07343"),
    @("E5", "-2.76%"),
    @("D6", "2.174"),
    @("E6", "43.01%"),
    @("D7", "7.927"),
    @("E7", "-0.01%"),
    @("D8", "3.791"),
    @("E8", "-0.81%"),
    @("D9", "0.9189"),
    @("E9", "-0.49%"),
    @("D10", "0.1702"),
    @("E10", "0.52%"),
    @("D11", "0.07528"),
    @("E11", "-4.79%"),
    @("D12", "0.08166"),
    @("E12", "1.41%"),
    @("D13", "0.03032"),
    @("E13", "0.07%"),
    @("D14", "0.09937"),
    @("E14", "0.26%"),
    @("D15", "0.001497"),
    @("E15", "-0.45%"),
    @("D16", "0.006120"),
    @("E16", "-3.26%"),
    @("E17", "0.42%"),
    @("D18", "2.220"),
    @("E18", "-0.71%"),
    @("D19", "0.3281"),
    @("E19", "-0.55%"),
    @("E20", "0.96%"),
    @("D21", "4.650"),
    @("E21", "3.92%"),
    @("D22", "0.04636"),
    @("E22", "0.90%"),
    @("D23", "0.1570"),
    @("E23", "-2.97%"),
    @("D24", "0.001227"),
    @("E24", "0.84%"),
    @("D25", "0.004464"),
    @("E25", "0.06%"),
    @("D26", "0.0001300"),
    @("E26", "-7.04%"),
    @("E27", "49.60%"),
    @("D39", "0.01724"),
    @("E39", "1.60%"),
    @("D40", "0.04515"),
    @("E40", "0.64%"),
    @("D41", "0.007289"),
    @("E41", "4.56%"),
    @("D42", "0.1345"),
    @("E43", "7.32%"),
    @("D44", "0.01063"),
    @("E44", "-22.93%"),
    @("D45", "0.00006333"),
    @("E45", "2.53%"),
    @("E46", "-23.03%"),
    @("E47", "15.83%")
)

foreach ($pair in $changes) {
    $ref = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}
